# Lesson 3 - "a few additions to lesson 3"
#
# Adds a second slide (Title and Content layout) that points to the
# Dublin Core usage-guide reference used in the second class meeting.

$p = $ppt.ActivePresentation

# Layout index 2 on the slide master is "Title and Content".
$s = $p.Slides.Add($p.Slides.Count + 1, 2)

# Leave the title placeholder blank and fill the content placeholder
# with the reference link for the lesson.
$s.Shapes.Placeholders.Item(2).TextFrame.TextRange.Text = "http://dublincore.org/documents/usageguide/elements.shtml"
